$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.282299160957336
$ws.Range("B1").Value = 1.788605213165283
$ws.Range("C1").Value = 2.280906438827515
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 0.9588499069213867
